# Final todo application code
# - users sheet: append the subbareddy user row
# - todos sheet: rename "completed" header to "target_date", add a new
#   "completed" header in column G, and append the AI todo row

$wb = $excel.ActiveWorkbook

# --- users sheet -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("users")

# Force these as plain text so numeric/ID-looking values aren't coerced
# into numbers by the automatic type-detection Excel performs on Value.
$ws1.Range("A2:D2").NumberFormat = "@"

$ws1.Range("A2").Value = "subbareddy"
$ws1.Range("B2").Value = '$2a$10$xIlSWnVXpCLkzt0f1HYv2.K0yPcVwAMLBQCK.2xoNlfCIk5cHRHSy'
$ws1.Range("C2").Value = "subbareddyroyal@gmail.com"
$ws1.Range("D2").Value = "2026-01-30T14:36:29.610Z"

# --- todos sheet -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("todos")

# Column D was "completed"; it's now "target_date". A new "completed"
# column is appended at G.
$ws2.Range("D1").Value = "target_date"
$ws2.Range("G1").Value = "completed"

# Force text so the big numeric id and the yyyy-mm-dd date string are
# stored verbatim instead of being reinterpreted as a number/date serial.
$ws2.Range("A2:F2").NumberFormat = "@"

$ws2.Range("A2").Value = "1769785959313"
$ws2.Range("B2").Value = "Artificial intelligence"
$ws2.Range("C2").Value = "Artificial intelligence"
$ws2.Range("D2").Value = "2026-02-28"
$ws2.Range("E2").Value = "2026-01-30T15:12:39.313Z"
$ws2.Range("F2").Value = "2026-01-30T15:12:59.592Z"
$ws2.Range("G2").Value = $false
